# Refresh cryptocurrency prices (column D) and 1-hour volume-change
# percentages (column E) on the active sheet, per the day's data pull.
#
# Every data cell in D/E is stored as text in this workbook (several
# "prices" use a thousands-dot, e.g. "47.316.00", which is not a valid
# Excel number). Assigning a plain numeric-looking string such as
# "321.32" straight to .Value would make Excel silently reinterpret it
# as a real number, so for every D-column update we briefly force a
# Text ("@") number format before writing the value, then restore the
# default "Normal" style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.316.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.493.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.41%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0811"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.124"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.49%  "

$ws.Range("E14").Value = "  -1.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.882.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.499.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.847"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.231.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.86%  "

$ws.Range("E21").Value = "  +0.18%  "

$ws.Range("E22").Value = "  +11.83%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "245.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.22%  "

$ws.Range("E25").Value = "  +0.65%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("E27").Value = "  -1.64%  "

$ws.Range("E28").Value = "  +3.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.17%  "

$ws.Range("E30").Value = "  +2.37%  "

$ws.Range("E31").Value = "  -1.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.51%  "

$ws.Range("E33").Value = "  +4.95%  "

$ws.Range("E34").Value = "  -0.17%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.71"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.47%  "

$ws.Range("E38").Value = "  +2.28%  "

$ws.Range("E39").Value = "  -1.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.53%  "

$ws.Range("E42").Value = "  +0.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.44%  "

$ws.Range("E44").Value = "  +0.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.994.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.17%  "

$ws.Range("E46").Value = "  +2.07%  "

$ws.Range("E47").Value = "  -5.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.55%  "

$ws.Range("E50").Value = "  -4.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.60%  "

